# Applies the "Add files via upload" diff: collapses runs that were
# split apart by proofErr (spell/grammar) markers back into single
# runs with the same visible text, and relocates the _GoBack bookmark.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $result = $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $result) {
        throw "Find/Replace failed for: $old"
    }
}

# 1. "Task 0 : Explain what you are doing/ going to accomplish"
Replace-Text "Task 0 : Explain what you are doing/ going to accomplish" "Task 0 : Explain what you are doing/ going to accomplish"

# 2. "I will need a button ... more then 0 items left in stock for them to buy"
Replace-Text "I will need a button to be displayed to the user when there is more then 0 items left in stock for them to buy" "I will need a button to be displayed to the user when there is more then 0 items left in stock for them to buy"

# 3. "Def purchase_success:  This function will -1 from stock when user clicks to purchase something"
Replace-Text "Def purchase_success:  This function will -1 from stock when user clicks to purchase something" "Def purchase_success:  This function will -1 from stock when user clicks to purchase something"

# 4. "AT ROUTE '/purchase-succcess/ <food.id>"
Replace-Text "AT ROUTE ‘/purchase-succcess/ <food.id>" "AT ROUTE ‘/purchase-succcess/ <food.id>"

# 5. "PROGRAM purchase_success (item.id)"
Replace-Text "PROGRAM purchase_success (item.id)" "PROGRAM purchase_success (item.id)"

# "      SET item.id TO int(item.id)"
Replace-Text "      SET item.id TO int(item.id)" "      SET item.id TO int(item.id)"

# "      SET found_item TO NONE"
Replace-Text "      SET found_item TO NONE" "      SET found_item TO NONE"

# "            IF item.id EQUALS item_id"
Replace-Text "            IF item.id EQUALS item_id" "            IF item.id EQUALS item_id"

# "                  SET found_item TO item"
Replace-Text "                  SET found_item TO item" "                  SET found_item TO item"

# "      SET data TO dict(item EQUALS found_item)"
Replace-Text "      SET data TO dict(item EQUALS found_item)" "      SET data TO dict(item EQUALS found_item)"

# "      SET found_item._stock TO -= 1"  -> becomes two runs around a
# relocated _GoBack bookmark: "      SET found_item._stock TO minus" + " 1"
Replace-Text "      SET found_item._stock TO -= 1" "      SET found_item._stock TO minus 1"

# 6. "Task 14 : Evaluation"
Replace-Text "Task 14 : Evaluation" "Task 14 : Evaluation"

# 7. "How did your version turn out"
Replace-Text "How did your version turn out" "How did your version turn out"

# 8. Move the _GoBack bookmark from the final paragraph to the
#    "minus"/" 1" split point in the purchase_success pseudocode.
$d.Bookmarks("_GoBack").Delete()

$found = $d.Content.Find.Execute("SET found_item._stock TO minus 1", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the minus-1 line to re-anchor _GoBack"
}
$anchorRange = $d.Content.Duplicate
$anchorRange.Start = $d.Content.Start + ($d.Content.Text.IndexOf("TO minus") + ("TO minus").Length)
$anchorRange.End = $anchorRange.Start
$d.Bookmarks.Add("_GoBack", $anchorRange) | Out-Null
